$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Name" text in B3 (shared string changes from the old net name
# to the new accumulator register name).
$ws.Range("B3").Value = "firConvolutionLoopUnrollingFactor2_IP/U0/accumulator_reg_140"

# Row 2 (summary row) - Signal Rate value update
$ws.Range("A2").Value = 0.0000031468659926758846268

# Row 3 (detail row) - updated power metrics
$ws.Range("A3").Value = 0.0000030680600957566639408
$ws.Range("C3").Value = 3.582089900970459
$ws.Range("D3").Value = 1.7910449504852295
$ws.Range("E3").Value = 37.0
$ws.Range("F3").Value = 11.0

# Row 4 (detail row) - updated power metrics
$ws.Range("A4").Value = 0.0000000788059679734942620
$ws.Range("C4").Value = 0.1492539942264557
$ws.Range("D4").Value = 3.7313430309295654
$ws.Range("E4").Value = 7.0
$ws.Range("F4").Value = 2.0
